# The only meaningful content change in the target revision is the header
# text in cell A1 of Sheet1: "name" -> "names" (the rest of the diff is
# just Excel-version resave metadata: fileVersion/build numbers, window
# geometry, revision/xr GUIDs, and namespace bumps, none of which are
# real spreadsheet edits).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("A1").Value = "names"
